$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "page_23"
$ws.Range("A1:Z1").Merge()
$ws.Range("A1").Value = " MASTER PACKAGE"
